# fixed time problem and wrote something on sop
#
# Adds a new "extra measurements" header/value/unit block in row 2/3
# (columns J..AG) to the spring-constant workbook: bar length/weight,
# weight weight, camera angle, weight height/r, disc r/weight,
# ball weight/diameter, ring weight/diameter.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 headers + first-seen units (order matters: it controls the
#     shared-string table build order so indices line up with the
#     canonical file) -----------------------------------------------------
$ws.Range("J2").Value = "bar længde"
$ws.Range("S3").Value = "cm"
$ws.Range("L2").Value = "bar vægt"
$ws.Range("M3").Value = "kg"
$ws.Range("K3").Value = "m"
$ws.Range("N2").Value = "weight weight"
$ws.Range("Y3").Value = "g"
$ws.Range("P2").Value = "camera angle"
$ws.Range("Q3").Value = "deg"
$ws.Range("R2").Value = "weight height"
$ws.Range("T2").Value = "weight r"
$ws.Range("V2").Value = "disc r"
$ws.Range("X2").Value = "disc weight"
$ws.Range("Z2").Value = "ball weight"
$ws.Range("AB2").Value = "ball diameter"
$ws.Range("AD2").Value = "ring weight"
$ws.Range("AF2").Value = "ring dia"

# --- Remaining repeated unit cells (reuse existing shared strings) ------
$ws.Range("U3").Value = "cm"
$ws.Range("W3").Value = "cm"
$ws.Range("AC3").Value = "cm"
$ws.Range("AG3").Value = "cm"
$ws.Range("O3").Value = "kg"
$ws.Range("AA3").Value = "g"
$ws.Range("AE3").Value = "g"

# --- Row 3 numeric measurements ------------------------------------------
$ws.Range("J3").Value = 0.6
$ws.Range("L3").Value = 0.128
$ws.Range("N3").Value = 0.227
$ws.Range("P3").Value = 22
$ws.Range("R3").Value = 4
$ws.Range("T3").Value = 1.5
$ws.Range("V3").Value = 10.8
$ws.Range("X3").Value = 270
$ws.Range("Z3").Value = 860
$ws.Range("AB3").Value = 14
$ws.Range("AD3").Value = 351
$ws.Range("AF3").Value = 10

# --- Scroll / selection, matching the author's recorded view state ------
$ws.Range("AH2").Select()
